$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 396.25
$ws.Range("I12").Value = 395
$ws.Range("J12").Value = 400
$ws.Range("K12").Value = 395
$ws.Range("L12").Value = 400
$ws.Range("M12").Value = -225
$ws.Range("N12").Value = -740

$ws.Range("H53").Value = 408.30768
$ws.Range("I53").Value = 95.833336
$ws.Range("K53").Value = 95.833336
$ws.Range("M53").Value = 541.166664

$ws.Range("H113").Value = 8583
$ws.Range("I113").Value = 7526.091
$ws.Range("J113").Value = 9413.429
$ws.Range("K113").Value = 7526.091
$ws.Range("L113").Value = 9413.429
$ws.Range("M113").Value = -4272.091
$ws.Range("N113").Value = -15921.429

$ws.Range("H130").Value = 28833
$ws.Range("J130").Value = 28833
$ws.Range("L130").Value = 28833
$ws.Range("N130").Value = -38873

$ws.Range("H141").Value = 5475.5
$ws.Range("I141").Value = 4382.364
$ws.Range("J141").Value = 17500
$ws.Range("K141").Value = 13147.092
$ws.Range("L141").Value = 52500
$ws.Range("M141").Value = -7967.091999999999
$ws.Range("N141").Value = -62860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14707203
$ws.Range("I32").Value = 15626280
$ws.Range("K32").Value = 15626280
$ws.Range("M32").Value = -15625993

$ws.Range("H132").Value = 2315.4517
$ws.Range("I132").Value = 1785.5172
$ws.Range("K132").Value = 5356.5516
$ws.Range("M132").Value = -2826.5516

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2409.375
$ws.Range("I20").Value = 1479
$ws.Range("K20").Value = 1479
$ws.Range("M20").Value = -1232

$ws.Range("H76").Value = 24790.5
$ws.Range("J76").Value = 24790.5
$ws.Range("L76").Value = 24790.5
$ws.Range("N76").Value = -25420.5

$ws.Range("H79").Value = 24790.5
$ws.Range("J79").Value = 24790.5
$ws.Range("L79").Value = 24790.5
$ws.Range("N79").Value = -26974.5

$ws.Range("H86").Value = 3260.8667
$ws.Range("I86").Value = 2295.5334
$ws.Range("K86").Value = 2295.5334
$ws.Range("M86").Value = -1172.5334

$ws.Range("H88").Value = 30874.75
$ws.Range("J88").Value = 30874.75
$ws.Range("L88").Value = 30874.75
$ws.Range("N88").Value = -31686.75

$ws.Range("H89").Value = 3260.8667
$ws.Range("I89").Value = 2295.5334
$ws.Range("K89").Value = 11477.667
$ws.Range("M89").Value = -5861.666999999999

$ws.Range("H91").Value = 30874.75
$ws.Range("J91").Value = 30874.75
$ws.Range("L91").Value = 30874.75
$ws.Range("N91").Value = -33682.75

$ws.Range("H105").Value = 2004.091
$ws.Range("I105").Value = 1923.8889
$ws.Range("K105").Value = 1923.8889
$ws.Range("M105").Value = -176.8888999999999

$ws.Range("H107").Value = 14027
$ws.Range("J107").Value = 23337.334
$ws.Range("L107").Value = 23337.334
$ws.Range("N107").Value = -27177.334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1432.909
$ws.Range("I31").Value = 1432.909
$ws.Range("K31").Value = 1432.909
$ws.Range("M31").Value = -1137.909

$ws.Range("H34").Value = 1432.909
$ws.Range("I34").Value = 1432.909
$ws.Range("K34").Value = 1432.909
$ws.Range("M34").Value = -1230.909

$ws.Range("H58").Value = 1648.2333
$ws.Range("I58").Value = 1286.1305
$ws.Range("J58").Value = 2838
$ws.Range("K58").Value = 1286.1305
$ws.Range("L58").Value = 2838
$ws.Range("M58").Value = -1083.1305
$ws.Range("N58").Value = -3244

$ws.Range("H86").Value = 27368.055
$ws.Range("I86").Value = 38985.555
$ws.Range("J86").Value = 15750.556
$ws.Range("K86").Value = 38985.555
$ws.Range("L86").Value = 15750.556
$ws.Range("M86").Value = -37862.555
$ws.Range("N86").Value = -17996.556

$ws.Range("H89").Value = 27368.055
$ws.Range("I89").Value = 38985.555
$ws.Range("J89").Value = 15750.556
$ws.Range("K89").Value = 194927.775
$ws.Range("L89").Value = 78752.78
$ws.Range("M89").Value = -189311.775
$ws.Range("N89").Value = -89984.78

$ws.Range("H132").Value = 2256.516
$ws.Range("I132").Value = 2044.3846
$ws.Range("K132").Value = 6133.1538
$ws.Range("M132").Value = -3603.1538

$ws.Range("H136").Value = 1648.2333
$ws.Range("I136").Value = 1286.1305
$ws.Range("J136").Value = 2838
$ws.Range("K136").Value = 3858.3915
$ws.Range("L136").Value = 8514
$ws.Range("M136").Value = -1308.3915
$ws.Range("N136").Value = -13614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 1506
$ws.Range("I76").Value = 1506
$ws.Range("K76").Value = 4518
$ws.Range("M76").Value = -4135

$ws.Range("H79").Value = 1506
$ws.Range("I79").Value = 1506
$ws.Range("K79").Value = 4518
$ws.Range("M79").Value = -3192

$ws.Range("H107").Value = 1142.5
$ws.Range("I107").Value = 434.5
$ws.Range("J107").Value = 1614.5
$ws.Range("K107").Value = 1303.5
$ws.Range("L107").Value = 4843.5
$ws.Range("M107").Value = 616.5
$ws.Range("N107").Value = -8683.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 3000
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 3000
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 3000
$ws.Range("M70").ClearContents()
$ws.Range("N70").Value = -3540

$ws.Range("H73").Value = 3000
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 3000
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 3000
$ws.Range("M73").ClearContents()
$ws.Range("N73").Value = -4872

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 593.6667
$ws.Range("I55").Value = 1011
$ws.Range("J55").Value = 259.8
$ws.Range("K55").Value = 1011
$ws.Range("L55").Value = 259.8
$ws.Range("M55").Value = -838
$ws.Range("N55").Value = -605.8

$ws.Range("H68").Value = 2195.3076
$ws.Range("I68").Value = 2240.0908
$ws.Range("J68").Value = 1949
$ws.Range("K68").Value = 2240.0908
$ws.Range("L68").Value = 1949
$ws.Range("M68").Value = -1491.0908
$ws.Range("N68").Value = -3447

$ws.Range("H71").Value = 2195.3076
$ws.Range("I71").Value = 2240.0908
$ws.Range("J71").Value = 1949
$ws.Range("K71").Value = 11200.454
$ws.Range("L71").Value = 9745
$ws.Range("M71").Value = -7456.454
$ws.Range("N71").Value = -17233

$ws.Range("H100").Value = 6797.6665
$ws.Range("I100").Value = 4082.3333
$ws.Range("K100").Value = 4082.3333
$ws.Range("M100").Value = -3541.3333

$ws.Range("H120").Value = 50232.668
$ws.Range("J120").Value = 50232.668
$ws.Range("L120").Value = 50232.668
$ws.Range("N120").Value = -59908.668

$ws.Range("H140").Value = 48424.25
$ws.Range("I140").Value = 30234
$ws.Range("J140").Value = 66614.5
$ws.Range("K140").Value = 30234
$ws.Range("L140").Value = 66614.5
$ws.Range("M140").Value = -25054
$ws.Range("N140").Value = -76974.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
